$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSMI_Map")

# Remove the duplicate "Specific.Conductance" row (row 39) from CSMI_Map;
# the "SpecificConductance" row (row 38) is kept.
$ws.Rows.Item(39).Delete()

# Keep the sheet's _FilterDatabase defined name range in sync with the
# reduced row count (was $H$1:$H$75, now one row shorter).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "CSMI_Map!_FilterDatabase") {
        $n.RefersTo = '=CSMI_Map!$H$1:$H$74'
    }
}
